# Auto-generated script to update currentAveragePrice / Leve price & profit columns
# per scheduled market-data refresh (columns H-N) across several Leve tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(13, 8).Value = 3844
$ws.Cells.Item(13, 9).Value = 305
$ws.Cells.Item(13, 10).Value = 18000
$ws.Cells.Item(13, 11).Value = 305
$ws.Cells.Item(13, 12).Value = 18000
$ws.Cells.Item(13, 13).Value = -136
$ws.Cells.Item(13, 14).Value = -18338

$ws.Cells.Item(19, 8).Value = 872.0769
$ws.Cells.Item(19, 9).Value = 1600.2
$ws.Cells.Item(19, 10).Value = 417
$ws.Cells.Item(19, 11).Value = 1600.2
$ws.Cells.Item(19, 12).Value = 417
$ws.Cells.Item(19, 13).Value = -1425.2
$ws.Cells.Item(19, 14).Value = -767

$ws.Cells.Item(41, 8).Value = 1383.1666
$ws.Cells.Item(41, 9).Value = 1972
$ws.Cells.Item(41, 10).Value = 558.8
$ws.Cells.Item(41, 11).Value = 1972
$ws.Cells.Item(41, 12).Value = 558.8
$ws.Cells.Item(41, 13).Value = -1532
$ws.Cells.Item(41, 14).Value = -1438.8

$ws.Cells.Item(42, 8).Value = 266.53845
$ws.Cells.Item(42, 9).Value = 206.875
$ws.Cells.Item(42, 10).Value = 362
$ws.Cells.Item(42, 11).Value = 620.625
$ws.Cells.Item(42, 12).Value = 1086
$ws.Cells.Item(42, 13).Value = -390.625

$ws.Cells.Item(43, 8).Value = 1373
$ws.Cells.Item(43, 9).Value = 985.5
$ws.Cells.Item(43, 10).Value = 1566.75
$ws.Cells.Item(43, 11).Value = 985.5
$ws.Cells.Item(43, 12).Value = 1566.75
$ws.Cells.Item(43, 13).Value = -916.5
$ws.Cells.Item(43, 14).Value = -1704.75

$ws.Cells.Item(51, 8).Value = 8549483
$ws.Cells.Item(51, 9).Value = 3950
$ws.Cells.Item(51, 10).Value = 10103216
$ws.Cells.Item(51, 11).Value = 3950
$ws.Cells.Item(51, 12).Value = 10103216
$ws.Cells.Item(51, 13).Value = -3466
$ws.Cells.Item(51, 14).Value = -10104184

$ws.Cells.Item(137, 8).Value = 22899.543
$ws.Cells.Item(137, 9).Value = 1034.9354
$ws.Cells.Item(137, 10).Value = 68086.39999999999
$ws.Cells.Item(137, 11).Value = 3104.8062
$ws.Cells.Item(137, 12).Value = 204259.2
$ws.Cells.Item(137, 13).Value = -554.8062
$ws.Cells.Item(137, 14).Value = -209359.2

$ws.Cells.Item(138, 8).Value = 3480.7307
$ws.Cells.Item(138, 9).Value = 2976.5557
$ws.Cells.Item(138, 10).Value = 3586.2559
$ws.Cells.Item(138, 11).Value = 8929.667099999999
$ws.Cells.Item(138, 12).Value = 10758.7677
$ws.Cells.Item(138, 13).Value = -3789.667099999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10189.21
$ws.Cells.Item(32, 9).Value = 5263.1665
$ws.Cells.Item(32, 10).Value = 27654.273
$ws.Cells.Item(32, 11).Value = 5263.1665
$ws.Cells.Item(32, 12).Value = 27654.273
$ws.Cells.Item(32, 13).Value = -4976.1665
$ws.Cells.Item(32, 14).Value = -28228.273

$ws.Cells.Item(74, 8).Value = 2656.5854
$ws.Cells.Item(74, 9).Value = 2631.3572
$ws.Cells.Item(74, 10).Value = 2710.923
$ws.Cells.Item(74, 11).Value = 2631.3572
$ws.Cells.Item(74, 12).Value = 2710.923
$ws.Cells.Item(74, 13).Value = -1757.3572
$ws.Cells.Item(74, 14).Value = -4458.923

$ws.Cells.Item(77, 8).Value = 2656.5854
$ws.Cells.Item(77, 9).Value = 2631.3572
$ws.Cells.Item(77, 10).Value = 2710.923
$ws.Cells.Item(77, 11).Value = 13156.786
$ws.Cells.Item(77, 12).Value = 13554.615
$ws.Cells.Item(77, 13).Value = -8788.786
$ws.Cells.Item(77, 14).Value = -22290.615

$ws.Cells.Item(102, 8).Value = 1648.4445
$ws.Cells.Item(102, 9).Value = 1648.4445
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 1648.4445
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = -26.44450000000006

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 582.5
$ws.Cells.Item(64, 9).Value = 306
$ws.Cells.Item(64, 10).Value = 613.2222
$ws.Cells.Item(64, 11).Value = 306
$ws.Cells.Item(64, 12).Value = 613.2222
$ws.Cells.Item(64, 13).Value = -81
$ws.Cells.Item(64, 14).Value = -1063.2222

$ws.Cells.Item(67, 8).Value = 582.5
$ws.Cells.Item(67, 9).Value = 306
$ws.Cells.Item(67, 10).Value = 613.2222
$ws.Cells.Item(67, 11).Value = 306
$ws.Cells.Item(67, 12).Value = 613.2222
$ws.Cells.Item(67, 13).Value = 474
$ws.Cells.Item(67, 14).Value = -2173.2222

$ws.Cells.Item(94, 8).Value = 9775.305
$ws.Cells.Item(94, 9).Value = 953.3333
$ws.Cells.Item(94, 10).Value = 26316.5
$ws.Cells.Item(94, 11).Value = 953.3333
$ws.Cells.Item(94, 12).Value = 26316.5
$ws.Cells.Item(94, 13).Value = -502.3333
$ws.Cells.Item(94, 14).Value = -27218.5

$ws.Cells.Item(105, 8).Value = 1693.1428
$ws.Cells.Item(105, 9).Value = 903
$ws.Cells.Item(105, 10).Value = 2746.6667
$ws.Cells.Item(105, 11).Value = 903
$ws.Cells.Item(105, 12).Value = 2746.6667
$ws.Cells.Item(105, 13).Value = 844
$ws.Cells.Item(105, 14).Value = -6240.6667

$ws.Cells.Item(134, 8).Value = 346673.75
$ws.Cells.Item(134, 9).Value = 556852.2
$ws.Cells.Item(134, 10).Value = 2745.4546
$ws.Cells.Item(134, 11).Value = 1670556.6
$ws.Cells.Item(134, 12).Value = 8236.363799999999
$ws.Cells.Item(134, 13).Value = -1668021.6
$ws.Cells.Item(134, 14).Value = -13306.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 1957.3334
$ws.Cells.Item(132, 9).Value = 1909.2941
$ws.Cells.Item(132, 10).Value = 2074
$ws.Cells.Item(132, 11).Value = 5727.8823
$ws.Cells.Item(132, 12).Value = 6222
$ws.Cells.Item(132, 13).Value = -3197.8823
$ws.Cells.Item(132, 14).Value = -11282

$ws.Cells.Item(134, 8).Value = 3006.0815
$ws.Cells.Item(134, 9).Value = 3109.8948
$ws.Cells.Item(134, 10).Value = 2647.4546
$ws.Cells.Item(134, 11).Value = 9329.6844
$ws.Cells.Item(134, 12).Value = 7942.3638
$ws.Cells.Item(134, 13).Value = -6794.6844

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(10, 8).Value = 12.5
$ws.Cells.Item(10, 9).Value = 12.5
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 37.5
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).Value = 101.5

$ws.Cells.Item(41, 8).Value = 57.77778
$ws.Cells.Item(41, 9).Value = 50
$ws.Cells.Item(41, 10).Value = 120
$ws.Cells.Item(41, 11).Value = 150
$ws.Cells.Item(41, 12).Value = 360
$ws.Cells.Item(41, 13).Value = 188
$ws.Cells.Item(41, 14).Value = -1036

$ws.Cells.Item(42, 8).Value = 8750
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 8750
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 12).Value = 26250
$ws.Cells.Item(42, 14).Value = -27318

$ws.Cells.Item(43, 8).Value = 6000
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 6000
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 18000
$ws.Cells.Item(43, 14).Value = -18228

$ws.Cells.Item(131, 8).Value = 776.4545000000001
$ws.Cells.Item(131, 9).Value = 358.07144
$ws.Cells.Item(131, 10).Value = 1084.7368
$ws.Cells.Item(131, 11).Value = 1074.21432
$ws.Cells.Item(131, 12).Value = 3254.2104
$ws.Cells.Item(131, 13).Value = 3965.78568
$ws.Cells.Item(131, 14).Value = -13334.2104

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 3214.8
$ws.Cells.Item(61, 9).Value = 3148
$ws.Cells.Item(61, 10).Value = 3482
$ws.Cells.Item(61, 11).Value = 3148
$ws.Cells.Item(61, 12).Value = 3482
$ws.Cells.Item(61, 13).Value = -2946
$ws.Cells.Item(61, 14).Value = -3886

$ws.Cells.Item(113, 8).Value = 3214.8
$ws.Cells.Item(113, 9).Value = 3148
$ws.Cells.Item(113, 10).Value = 3482
$ws.Cells.Item(113, 11).Value = 3148
$ws.Cells.Item(113, 12).Value = 3482
$ws.Cells.Item(113, 13).Value = -978
$ws.Cells.Item(113, 14).Value = -7822

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(9, 8).Value = 795
$ws.Cells.Item(9, 9).Value = 795
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 795
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 13).Value = -655
$ws.Cells.Item(9, 14).ClearContents()
